# Applies the Mon Apr 17 22:38:01 UTC 2023 "cryptos" price/volume refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '29.740.74'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -2.49%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.095.97'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -1.81%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.010'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '344.02'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -2.33%  '

$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("E7").Value = '  -1.65%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.4381'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -3.72%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '52.59'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -2.03%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.09267'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +1.44%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '1.163'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -2.46%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '24.89'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -2.12%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '2.104.25'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -1.86%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '8.259'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +1.34%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '6.753'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -1.81%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '99.54'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -2.01%  '

$ws.Range("E17").Value = '  -1.08%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '1.009'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +0.19%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '20.81'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +1.18%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.06650'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -1.00%  '

$ws.Range("E21").Value = '  +0.28%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '6.192'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -2.72%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '29.763.52'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -2.73%  '

$ws.Range("E24").Value = '  -2.96%  '

$ws.Range("E25").Value = '  -2.88%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.348.92'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -1.87%  '

$ws.Range("E27").Value = '  -2.51%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '2.510'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -4.02%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '161.28'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -2.20%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '133.01'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -1.97%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '1.142'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -6.49%  '

$ws.Range("E32").Value = '  -3.10%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.650'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -3.70%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '6.170'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -3.36%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '3.936'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -2.30%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '6.284'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +2.73%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '10.22'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -2.00%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.02580'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -2.55%  '

$ws.Range("E39").Value = '  -3.53%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '12.46'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -2.07%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.6886'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -1.33%  '

$ws.Range("E42").Value = '  -5.20%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '1.317'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +3.52%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.6769'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +3.97%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '14.28'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -3.50%  '

$ws.Range("E46").Value = '  -1.13%  '

$ws.Range("E47").Value = '  -4.85%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '3.621'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -3.35%  '

$ws.Range("E49").Value = '  -2.16%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '82.07'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -1.95%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.160'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -2.14%  '
